$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "A burndown chart I created to record the teams sprint progress "
$ws.Range("A13").Value = "throughout Release 1"

$ws.Range("B16").Select()
